$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "party_2020" column header in AC1
$ws.Range("AC1").Value = "party_2020"

# Populate party_2020 for each state row (2020 US presidential election winner per state)
$ws.Range("AC2").Value = "REPUBLICAN"
$ws.Range("AC3").Value = "REPUBLICAN"
$ws.Range("AC4").Value = "DEMOCRAT"
$ws.Range("AC5").Value = "REPUBLICAN"
$ws.Range("AC6").Value = "DEMOCRAT"
$ws.Range("AC7").Value = "DEMOCRAT"
$ws.Range("AC8").Value = "DEMOCRAT"
$ws.Range("AC9").Value = "DEMOCRAT"
$ws.Range("AC10").Value = "DEMOCRAT"
$ws.Range("AC11").Value = "REPUBLICAN"
$ws.Range("AC12").Value = "DEMOCRAT"
$ws.Range("AC13").Value = "DEMOCRAT"
$ws.Range("AC14").Value = "REPUBLICAN"
$ws.Range("AC15").Value = "DEMOCRAT"
$ws.Range("AC16").Value = "REPUBLICAN"
$ws.Range("AC17").Value = "REPUBLICAN"
$ws.Range("AC18").Value = "REPUBLICAN"
$ws.Range("AC19").Value = "REPUBLICAN"
$ws.Range("AC20").Value = "REPUBLICAN"
$ws.Range("AC21").Value = "DEMOCRAT"
$ws.Range("AC22").Value = "DEMOCRAT"
$ws.Range("AC23").Value = "DEMOCRAT"
$ws.Range("AC24").Value = "DEMOCRAT"
$ws.Range("AC25").Value = "DEMOCRAT"
$ws.Range("AC26").Value = "REPUBLICAN"
$ws.Range("AC27").Value = "REPUBLICAN"
$ws.Range("AC28").Value = "REPUBLICAN"
$ws.Range("AC29").Value = "REPUBLICAN"
$ws.Range("AC30").Value = "DEMOCRAT"
$ws.Range("AC31").Value = "DEMOCRAT"
$ws.Range("AC32").Value = "DEMOCRAT"
$ws.Range("AC33").Value = "DEMOCRAT"
$ws.Range("AC34").Value = "DEMOCRAT"
$ws.Range("AC35").Value = "REPUBLICAN"
$ws.Range("AC36").Value = "REPUBLICAN"
$ws.Range("AC37").Value = "REPUBLICAN"
$ws.Range("AC38").Value = "REPUBLICAN"
$ws.Range("AC39").Value = "DEMOCRAT"
$ws.Range("AC40").Value = "DEMOCRAT"
$ws.Range("AC41").Value = "DEMOCRAT"
$ws.Range("AC42").Value = "REPUBLICAN"
$ws.Range("AC43").Value = "REPUBLICAN"
$ws.Range("AC44").Value = "REPUBLICAN"
$ws.Range("AC45").Value = "REPUBLICAN"
$ws.Range("AC46").Value = "REPUBLICAN"
$ws.Range("AC47").Value = "DEMOCRAT"
$ws.Range("AC48").Value = "DEMOCRAT"
$ws.Range("AC49").Value = "DEMOCRAT"
$ws.Range("AC50").Value = "REPUBLICAN"
$ws.Range("AC51").Value = "DEMOCRAT"
$ws.Range("AC52").Value = "REPUBLICAN"

# Update selection to match the edited workbook's saved cursor position
$ws.Range("K18").Select() | Out-Null
